$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6 (pushes old rows 6..27 down to 7..28)
$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with the "Ignore test" content.
# The leading apostrophe forces "=" to be treated as literal text instead
# of being parsed as a formula; resetting the Style afterwards clears the
# "quote prefix" formatting flag that the apostrophe trick leaves behind.
$a6 = $ws.Cells.Item(6, 1)
$a6.Value = "'="
$a6.Style = "Normal"

$ws.Cells.Item(6, 2).Value = "Ignore test:"
$ws.Cells.Item(6, 3).Value = "xltablediff.py  --key ID --ignore Color test1old.xlsx test1new.xlsx --out test1ignore.xlsx"

Write-Output "Done"
